$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting for numeric-looking price strings
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.139.86'
$ws.Range("E2").Value = '  -0.08%  '

$ws.Range("D3").Value = '1.899.51'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").Value = '306.89'
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("D7").Value = '0.5234'
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  +0.78%  '

$ws.Range("D9").Value = '0.07283'
$ws.Range("E9").Value = '  +0.51%  '

$ws.Range("D10").Value = '21.32'
$ws.Range("E10").Value = '  +0.70%  '

$ws.Range("D11").Value = '0.9057'
$ws.Range("E11").Value = '  +0.79%  '

$ws.Range("D12").Value = '0.08201'
$ws.Range("E12").Value = '  -1.97%  '

$ws.Range("D13").Value = '1.891.43'
$ws.Range("E13").Value = '  -0.62%  '

$ws.Range("D14").Value = '95.63'
$ws.Range("E14").Value = '  +0.99%  '

$ws.Range("D15").Value = '5.351'
$ws.Range("E15").Value = '  +1.63%  '

$ws.Range("E16").Value = '  +0.24%  '

$ws.Range("D17").Value = '0.000008650'
$ws.Range("E17").Value = '  +0.78%  '

$ws.Range("E19").Value = '  +0.25%  '

$ws.Range("D20").Value = '27.182.88'
$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").Value = '5.120'
$ws.Range("E21").Value = '  +1.28%  '

$ws.Range("D22").Value = '2.123.14'
$ws.Range("E22").Value = '  -1.39%  '

$ws.Range("E23").Value = '  +1.97%  '

$ws.Range("D24").Value = '6.470'
$ws.Range("E24").Value = '  +0.79%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '149.51'
$ws.Range("E25").Value = '  +1.99%  '

$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.326'
$ws.Range("E26").Value = '  +1.91%  '

$ws.Range("D27").Value = '18.25'
$ws.Range("E27").Value = '  +0.84%  '

$ws.Range("D28").Value = '1.735'
$ws.Range("E28").Value = '  -1.15%  '

$ws.Range("D29").Value = '115.42'
$ws.Range("E29").Value = '  +0.67%  '

$ws.Range("D30").Value = '4.825'
$ws.Range("E30").Value = '  +0.99%  '

$ws.Range("D31").Value = '4.874'
$ws.Range("E31").Value = '  -0.96%  '

$ws.Range("E32").Value = '  +0.09%  '

$ws.Range("D33").Value = '0.05046'
$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("D34").Value = '0.7930'
$ws.Range("E34").Value = '  -3.19%  '

$ws.Range("D35").Value = '1.224'
$ws.Range("E35").Value = '  -0.90%  '

$ws.Range("D36").Value = '2.977'
$ws.Range("E36").Value = '  +0.65%  '

$ws.Range("D37").Value = '3.375'
$ws.Range("E37").Value = '  +0.31%  '

$ws.Range("D38").Value = '2.644'
$ws.Range("E38").Value = '  +2.86%  '

$ws.Range("D39").Value = '0.5739'
$ws.Range("E39").Value = '  +0.90%  '

$ws.Range("E40").Value = '  +1.11%  '

$ws.Range("D41").Value = '1.081'
$ws.Range("E41").Value = '  +0.67%  '

$ws.Range("D42").Value = '9.038'
$ws.Range("E42").Value = '  +1.15%  '

$ws.Range("D43").Value = '6.622'
$ws.Range("E43").Value = '  -0.47%  '

$ws.Range("D44").Value = '116.38'
$ws.Range("E44").Value = '  -1.56%  '

$ws.Range("E45").Value = '  +0.49%  '

$ws.Range("D46").Value = '0.4902'
$ws.Range("E46").Value = '  +1.77%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '10.19'
$ws.Range("E47").Value = '  +0.54%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '1.003'
$ws.Range("E48").Value = '  +0.23%  '

$ws.Range("E49").Value = '  +2.08%  '

$ws.Range("D50").Value = '38.59'
$ws.Range("E50").Value = '  +3.08%  '

$ws.Range("D51").Value = '64.13'
$ws.Range("E51").Value = '  +0.91%  '
